# Fruta / hortaliza, semanal
# Insert a new weekly record as row 14 (pushing the previous rows 14-20
# down to 15-21), matching the new price-report row appended upstream.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift rows 14:20 down to 15:21 by inserting a new blank row at 14.
$ws.Rows.Item(14).Insert()

# Populate the new row 14 with the latest weekly data point.
$ws.Cells.Item(14, 1).Value = 4
$ws.Cells.Item(14, 2).Value = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(14, 3).Value = "Los Lagos"
$ws.Cells.Item(14, 4).Value = 44572
$ws.Cells.Item(14, 5).Value = 10
$ws.Cells.Item(14, 6).Value = 100112030
$ws.Cells.Item(14, 7).Value = "Poroto granado"
$ws.Cells.Item(14, 8).Value = "Sin especificar"
$ws.Cells.Item(14, 9).Value = "Primera"
$ws.Cells.Item(14, 10).Value = 80
$ws.Cells.Item(14, 11).Value = 35000
$ws.Cells.Item(14, 12).Value = 35000
$ws.Cells.Item(14, 13).Value = 35000
$ws.Cells.Item(14, 14).Value = "$/saco 25 kilos"
$ws.Cells.Item(14, 15).Value = "Región Metropolitana"
$ws.Cells.Item(14, 16).Value = 1400
$ws.Cells.Item(14, 17).Value = 25
$ws.Cells.Item(14, 18).Value = "Hortaliza"
